$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "modif nicollini / ordinaire" -- flip the validation flag (column C) for
# niccolini_lezioni_1880.xml (row 26) and ordinaire_dictionnaire-mythologie_1866.xml (row 27)
$ws.Range("C26").Value = 1
$ws.Range("C27").Value = 1

# Scroll the view down a bit and leave the selection on C26, matching where
# the author ended up after making the edit above.
$excel.ActiveWindow.ScrollRow = 3
$ws.Range("C26").Select()
